$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: split the run containing " and Romig et al., " so the word
# "Romig" sits in its own run (matches the source after Word's proofing
# pass re-flagged "Romig" as a spelling exception around it). Visible
# text is unchanged; only the run boundaries move. Toggling a character
# property on just that word and then reverting it forces Word to keep
# the run split (Word does not silently re-merge runs once split).
# ---------------------------------------------------------------------
$find1 = $d.Content
$found1 = $find1.Find.Execute("and Romig et al", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $full1 = $find1.Text
    $offset = $full1.IndexOf("Romig")
    $wordStart = $find1.Start + $offset
    $wordEnd = $wordStart + 5
    $romigRange = $d.Range($wordStart, $wordEnd)
    $romigRange.Font.Bold = $true
    $romigRange.Font.Bold = $false
}

# ---------------------------------------------------------------------
# Change 2: populate the empty paragraph directly under the "Table of
# Contents" heading with the placeholder "[ToC]" text.
# ---------------------------------------------------------------------
$found2 = $false
foreach ($p in $d.Paragraphs) {
    if ($found2) {
        $p.Range.Text = "[ToC]"
        break
    }
    if ($p.Range.Text.TrimEnd("`r") -eq "Table of Contents") {
        $found2 = $true
    }
}

# ---------------------------------------------------------------------
# Change 3: the phrase "... School of [Law, and] published in 2021 ..."
# had "Law, and" split into its own run flanked by grammar-check markers;
# collapse it back into the single surrounding run (visible text is
# unchanged). Rewriting the whole span as one assignment merges it back
# into a single run.
# ---------------------------------------------------------------------
$find3 = $d.Content
$target3 = " was compiled by Professor Jennifer Murphy Romig with assistance from students at the Emory University School of Law, and published in 2021. In the five years since the first edition of "
$found3 = $find3.Find.Execute("was compiled by Professor Jennifer Murphy Romig with assistance", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $full3 = $find3.Text
    $spanStart = $find3.Start - 1
    $spanEnd = $spanStart + $target3.Length
    $spanRange = $d.Range($spanStart, $spanEnd)
    # Sanity-check the text before touching it.
    if ($spanRange.Text -eq $target3) {
        $spanRange.Text = "PLACEHOLDER_MERGE_TOKEN"
        $reinsertRange = $d.Range($spanStart, $spanStart + 23)
        $reinsertRange.Text = $target3
    }
}
